$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 201.33333
$ws.Range("I9").Value = 305.5
$ws.Range("J9").Value = 118
$ws.Range("K9").Value = 305.5
$ws.Range("L9").Value = 118
$ws.Range("M9").Value = -136.5
$ws.Range("N9").Value = -456

$ws.Range("H15").Value = 2377.1233
$ws.Range("I15").Value = 2377.1233
$ws.Range("K15").Value = 7131.369900000001
$ws.Range("M15").Value = -6962.369900000001

$ws.Range("H17").Value = 515024.97
$ws.Range("J17").Value = 515024.97
$ws.Range("L17").Value = 1545074.91
$ws.Range("N17").Value = -1545410.91

$ws.Range("H19").Value = 601.875
$ws.Range("I19").Value = 514.5
$ws.Range("J19").Value = 689.25
$ws.Range("K19").Value = 514.5
$ws.Range("L19").Value = 689.25
$ws.Range("M19").Value = -339.5
$ws.Range("N19").Value = -1039.25

$ws.Range("H38").Value = 533.3
$ws.Range("I38").Value = 15.75
$ws.Range("J38").Value = 878.3333
$ws.Range("K38").Value = 47.25
$ws.Range("L38").Value = 2634.9999
$ws.Range("M38").Value = 324.75
$ws.Range("N38").Value = -3378.9999

$ws.Range("H42").Value = 227.5
$ws.Range("I42").Value = 55
$ws.Range("J42").Value = 400
$ws.Range("K42").Value = 165
$ws.Range("L42").Value = 1200
$ws.Range("M42").Value = 65
$ws.Range("N42").Value = -1660

$ws.Range("H64").Value = 83337090
$ws.Range("I64").Value = 250001870
$ws.Range("J64").Value = 4695
$ws.Range("K64").Value = 250001870
$ws.Range("L64").Value = 4695
$ws.Range("M64").Value = -250001622
$ws.Range("N64").Value = -5191

$ws.Range("H67").Value = 83337090
$ws.Range("I67").Value = 250001870
$ws.Range("J67").Value = 4695
$ws.Range("K67").Value = 250001870
$ws.Range("L67").Value = 4695
$ws.Range("M67").Value = -250001012
$ws.Range("N67").Value = -6411

$ws.Range("H74").Value = 3460.8667
$ws.Range("I74").Value = 2756.7693
$ws.Range("J74").Value = 3999.2942
$ws.Range("K74").Value = 2756.7693
$ws.Range("L74").Value = 3999.2942
$ws.Range("M74").Value = -1820.7693
$ws.Range("N74").Value = -5871.2942

$ws.Range("H77").Value = 3460.8667
$ws.Range("I77").Value = 2756.7693
$ws.Range("J77").Value = 3999.2942
$ws.Range("K77").Value = 13783.8465
$ws.Range("L77").Value = 19996.471
$ws.Range("M77").Value = -9103.8465
$ws.Range("N77").Value = -29356.471

$ws.Range("H92").Value = 282.8
$ws.Range("I92").Value = 233.66667
$ws.Range("J92").Value = 725
$ws.Range("K92").Value = 233.66667
$ws.Range("L92").Value = 725
$ws.Range("M92").Value = 1014.33333
$ws.Range("N92").Value = -3221

$ws.Range("H135").Value = 579.9643
$ws.Range("I135").Value = 567.4167
$ws.Range("J135").Value = 655.25
$ws.Range("K135").Value = 5106.7503
$ws.Range("L135").Value = 5897.25
$ws.Range("M135").Value = -2571.7503
$ws.Range("N135").Value = -10967.25

$ws.Range("H138").Value = 2354.7036
$ws.Range("I138").Value = 1017.7143
$ws.Range("J138").Value = 3794.5386
$ws.Range("K138").Value = 3053.1429
$ws.Range("L138").Value = 11383.6158
$ws.Range("M138").Value = 2086.8571
$ws.Range("N138").Value = -21663.6158

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 883.89
$ws.Range("I32").Value = 676.6484
$ws.Range("J32").Value = 2979.3333
$ws.Range("K32").Value = 676.6484
$ws.Range("L32").Value = 2979.3333
$ws.Range("M32").Value = -389.6484
$ws.Range("N32").Value = -3553.3333

$ws.Range("H61").Value = 1040.0405
$ws.Range("I61").Value = 784.0339
$ws.Range("J61").Value = 2047
$ws.Range("K61").Value = 784.0339
$ws.Range("L61").Value = 2047
$ws.Range("M61").Value = -572.0339
$ws.Range("N61").Value = -2471

$ws.Range("H74").Value = 3028.9556
$ws.Range("I74").Value = 714.025
$ws.Range("J74").Value = 21548.4
$ws.Range("K74").Value = 714.025
$ws.Range("L74").Value = 21548.4
$ws.Range("M74").Value = 159.975
$ws.Range("N74").Value = -23296.4

$ws.Range("H77").Value = 3028.9556
$ws.Range("I77").Value = 714.025
$ws.Range("J77").Value = 21548.4
$ws.Range("K77").Value = 3570.125
$ws.Range("L77").Value = 107742
$ws.Range("M77").Value = 797.875
$ws.Range("N77").Value = -116478

$ws.Range("H136").Value = 1040.0405
$ws.Range("I136").Value = 784.0339
$ws.Range("J136").Value = 2047
$ws.Range("K136").Value = 2352.1017
$ws.Range("L136").Value = 6141
$ws.Range("M136").Value = 197.8982999999998
$ws.Range("N136").Value = -11241

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 992.8
$ws.Range("I94").Value = 653.1667
$ws.Range("J94").Value = 1866.1428
$ws.Range("K94").Value = 653.1667
$ws.Range("L94").Value = 1866.1428
$ws.Range("M94").Value = -202.1667
$ws.Range("N94").Value = -2768.1428

$ws.Range("H134").Value = 965.8788
$ws.Range("I134").Value = 888.86664
$ws.Range("J134").Value = 1736
$ws.Range("K134").Value = 2666.59992
$ws.Range("L134").Value = 5208
$ws.Range("M134").Value = -131.5999199999997
$ws.Range("N134").Value = -10278

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 35066.58
$ws.Range("I31").Value = 1099.3429
$ws.Range("J31").Value = 76061.516
$ws.Range("K31").Value = 1099.3429
$ws.Range("L31").Value = 76061.516
$ws.Range("M31").Value = -804.3429000000001
$ws.Range("N31").Value = -76651.516

$ws.Range("H34").Value = 35066.58
$ws.Range("I34").Value = 1099.3429
$ws.Range("J34").Value = 76061.516
$ws.Range("K34").Value = 1099.3429
$ws.Range("L34").Value = 76061.516
$ws.Range("M34").Value = -897.3429000000001
$ws.Range("N34").Value = -76465.516

$ws.Range("H37").Value = 16001.875
$ws.Range("I37").Value = 3836.5
$ws.Range("J37").Value = 20057
$ws.Range("K37").Value = 3836.5
$ws.Range("L37").Value = 20057
$ws.Range("M37").Value = -3729.5
$ws.Range("N37").Value = -20271

$ws.Range("H134").Value = 3665.2341
$ws.Range("I134").Value = 3884.244
$ws.Range("J134").Value = 2168.6667
$ws.Range("K134").Value = 11652.732
$ws.Range("L134").Value = 6506.000100000001
$ws.Range("M134").Value = -9117.732
$ws.Range("N134").Value = -11576.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 353105.25
$ws.Range("I113").Value = 423.32352
$ws.Range("J113").Value = 716474.5
$ws.Range("K113").Value = 1269.97056
$ws.Range("L113").Value = 2149423.5
$ws.Range("M113").Value = 900.02944
$ws.Range("N113").Value = -2153763.5

$ws.Range("H121").Value = 7320.7646
$ws.Range("I121").Value = 460
$ws.Range("J121").Value = 7984.7095
$ws.Range("K121").Value = 1380
$ws.Range("L121").Value = 23954.1285
$ws.Range("M121").Value = -70
$ws.Range("N121").Value = -26574.1285

$ws.Range("H131").Value = 777.7377300000001
$ws.Range("I131").Value = 571.4
$ws.Range("J131").Value = 977.4194
$ws.Range("K131").Value = 1714.2
$ws.Range("L131").Value = 2932.2582
$ws.Range("M131").Value = 3325.8
$ws.Range("N131").Value = -13012.2582

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5775.4165
$ws.Range("I132").Value = 8258.736999999999
$ws.Range("J132").Value = 2999.9412
$ws.Range("K132").Value = 24776.211
$ws.Range("L132").Value = 8999.8236
$ws.Range("M132").Value = -22246.211
$ws.Range("N132").Value = -14059.8236

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 249.76923
$ws.Range("I107").Value = 237.25
$ws.Range("J107").Value = 400
$ws.Range("K107").Value = 711.75
$ws.Range("L107").Value = 1200
$ws.Range("M107").Value = 1208.25
$ws.Range("N107").Value = -5040

$ws.Range("H132").Value = 2858.228
$ws.Range("I132").Value = 4329.968
$ws.Range("J132").Value = 1103.4615
$ws.Range("K132").Value = 12989.904
$ws.Range("L132").Value = 3310.3845
$ws.Range("M132").Value = -10459.904
$ws.Range("N132").Value = -8370.3845

$ws.Range("H136").Value = 2354.4028
$ws.Range("I136").Value = 2423
$ws.Range("J136").Value = 2208.261
$ws.Range("K136").Value = 7269
$ws.Range("L136").Value = 6624.782999999999
$ws.Range("M136").Value = -4719
$ws.Range("N136").Value = -11724.783
